# Update "想去人数" (attendance count) values in column F for the
# "展览" and "全部类型" worksheets, per the commit's regenerated data.

$wb = $excel.ActiveWorkbook

# row -> new value, for worksheet "展览" (sheet1)
$changesExhibition = @{
    2  = 1902
    4  = 869
    6  = 52
    8  = 253
    11 = 145
    13 = 4507
    14 = 16
    16 = 495
    17 = 448
    20 = 1221
    21 = 2372
    23 = 65
    24 = 44
    25 = 56
    26 = 2214
    30 = 159
    31 = 103
}

# row -> new value, for worksheet "全部类型" (sheet4)
$changesAllTypes = @{
    2  = 1902
    4  = 869
    6  = 52
    8  = 253
    11 = 145
    14 = 4507
    15 = 16
    17 = 495
    18 = 448
    21 = 1221
    22 = 2372
    24 = 65
    25 = 44
    26 = 56
    27 = 2214
    31 = 159
    32 = 103
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $changesExhibition.Keys) {
    $wsExhibition.Range("F$row").Value = $changesExhibition[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $changesAllTypes.Keys) {
    $wsAllTypes.Range("F$row").Value = $changesAllTypes[$row]
}
